$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.809.34'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").Value = '1.859.75'
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("D4").Value = '''1.037'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.69%  '

$ws.Range("D5").Value = '''323.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("E7").Value = '  +0.86%  '

$ws.Range("D8").Value = '''0.3819'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.70%  '

$ws.Range("D9").Value = '''0.07450'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").Value = '''0.8879'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.34%  '

$ws.Range("D11").Value = '''21.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.66%  '

$ws.Range("D12").Value = '1.863.01'
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("D13").Value = '''5.552'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("E14").Value = '  +0.67%  '

$ws.Range("D15").Value = '''0.07219'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("D16").Value = '''86.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.07%  '

$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '''0.000009118'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").Value = '''15.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").Value = '27.817.28'
$ws.Range("E21").Value = '  +0.65%  '

$ws.Range("D22").Value = '''5.303'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").Value = '''11.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '

$ws.Range("D24").Value = '2.077.90'
$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("E25").Value = '  +6.53%  '

$ws.Range("D26").Value = '''159.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.26%  '

$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("D28").Value = '''2.008'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.74%  '

$ws.Range("D29").Value = '''5.385'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").Value = '''118.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("D31").Value = '''0.09118'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").Value = '''1.219'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.90%  '

$ws.Range("D33").Value = '''0.7765'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("D34").Value = '''3.023'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.89%  '

$ws.Range("D35").Value = '''4.615'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.13%  '

$ws.Range("E36").Value = '  +0.56%  '

$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D39").Value = '''0.05325'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("D40").Value = '''2.863'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.07%  '

$ws.Range("D41").Value = '''0.5217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("D42").Value = '''6.980'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.54%  '

$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("D44").Value = '''8.815'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''111.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.06%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.73%  '

$ws.Range("D47").Value = '''1.037'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("D48").Value = '''0.06583'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.04%  '

$ws.Range("E49").Value = '  +0.13%  '

$ws.Range("D50").Value = '''0.4740'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '

$ws.Range("E51").Value = '  -0.14%  '

